$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 20010766
$ws.Cells.Item(92, 9).Value = 22233852
$ws.Cells.Item(92, 11).Value = 22233852
$ws.Cells.Item(92, 13).Value = -22232604
$ws.Cells.Item(94, 8).Value = 36743.4
$ws.Cells.Item(94, 9).Value = 36510.785
$ws.Cells.Item(94, 11).Value = 36510.785
$ws.Cells.Item(94, 13).Value = -36059.785
$ws.Cells.Item(106, 8).Value = 41668376
$ws.Cells.Item(106, 9).Value = 45456148
$ws.Cells.Item(106, 11).Value = 45456148
$ws.Cells.Item(106, 13).Value = -45455517
$ws.Cells.Item(112, 8).Value = 92567
$ws.Cells.Item(112, 9).Value = 653.5
$ws.Cells.Item(112, 10).Value = 145089
$ws.Cells.Item(112, 11).Value = 1960.5
$ws.Cells.Item(112, 12).Value = 435267
$ws.Cells.Item(112, 13).Value = -852.5
$ws.Cells.Item(112, 14).Value = -437483
$ws.Cells.Item(125, 8).Value = 642.75
$ws.Cells.Item(125, 9).Value = 642.75
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 11).Value = 5784.75
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 13).Value = -3324.75
$ws.Cells.Item(125, 14).ClearContents()
$ws.Cells.Item(135, 8).Value = 1555
$ws.Cells.Item(135, 9).Value = 1066
$ws.Cells.Item(135, 10).Value = 4000
$ws.Cells.Item(135, 11).Value = 9594
$ws.Cells.Item(135, 12).Value = 36000
$ws.Cells.Item(135, 13).Value = -7059
$ws.Cells.Item(135, 14).Value = -41070
$ws.Cells.Item(137, 8).Value = 2169
$ws.Cells.Item(137, 10).Value = 2000.25
$ws.Cells.Item(137, 12).Value = 6000.75
$ws.Cells.Item(137, 14).Value = -11100.75
$ws.Cells.Item(138, 8).Value = 3306.9048
$ws.Cells.Item(138, 9).Value = 2146.36
$ws.Cells.Item(138, 11).Value = 6439.08
$ws.Cells.Item(138, 13).Value = -1299.08

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3834.897
$ws.Cells.Item(32, 9).Value = 2936.0645
$ws.Cells.Item(32, 11).Value = 2936.0645
$ws.Cells.Item(32, 13).Value = -2649.0645
$ws.Cells.Item(74, 8).Value = 5930.921
$ws.Cells.Item(74, 9).Value = 1190.8788
$ws.Cells.Item(74, 11).Value = 1190.8788
$ws.Cells.Item(74, 13).Value = -316.8788
$ws.Cells.Item(77, 8).Value = 5930.921
$ws.Cells.Item(77, 9).Value = 1190.8788
$ws.Cells.Item(77, 11).Value = 5954.394
$ws.Cells.Item(77, 13).Value = -1586.394
$ws.Cells.Item(132, 8).Value = 3014.5557
$ws.Cells.Item(132, 9).Value = 2120.7058
$ws.Cells.Item(132, 10).Value = 4534.1
$ws.Cells.Item(132, 11).Value = 6362.117400000001
$ws.Cells.Item(132, 12).Value = 13602.3
$ws.Cells.Item(132, 13).Value = -3832.117400000001
$ws.Cells.Item(132, 14).Value = -18662.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1971.7142
$ws.Cells.Item(134, 9).Value = 1970.8823
$ws.Cells.Item(134, 11).Value = 5912.6469
$ws.Cells.Item(134, 13).Value = -3377.6469

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 13).ClearContents()
$ws.Cells.Item(59, 8).Value = 25470.588
$ws.Cells.Item(59, 10).Value = 25470.588
$ws.Cells.Item(59, 12).Value = 25470.588
$ws.Cells.Item(59, 14).Value = -27760.588
$ws.Cells.Item(88, 8).Value = 7900
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 14).ClearContents()
$ws.Cells.Item(91, 8).Value = 7900
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 855.8333
$ws.Cells.Item(107, 9).Value = 403.42856
$ws.Cells.Item(107, 11).Value = 403.42856
$ws.Cells.Item(107, 13).Value = 1516.57144
$ws.Cells.Item(131, 8).Value = 46333
$ws.Cells.Item(131, 9).Value = 38375
$ws.Cells.Item(131, 11).Value = 38375
$ws.Cells.Item(131, 13).Value = -33335
$ws.Cells.Item(132, 8).Value = 5011.222
$ws.Cells.Item(132, 9).Value = 4894.2354
$ws.Cells.Item(132, 11).Value = 14682.7062
$ws.Cells.Item(132, 13).Value = -12152.7062

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 501.5
$ws.Cells.Item(23, 10).Value = 501.5
$ws.Cells.Item(23, 12).Value = 1504.5
$ws.Cells.Item(23, 14).Value = -1974.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8347.962
$ws.Cells.Item(70, 9).Value = 4999.524
$ws.Cells.Item(70, 10).Value = 22411.4
$ws.Cells.Item(70, 11).Value = 4999.524
$ws.Cells.Item(70, 12).Value = 22411.4
$ws.Cells.Item(70, 13).Value = -4729.524
$ws.Cells.Item(70, 14).Value = -22951.4
$ws.Cells.Item(73, 8).Value = 8347.962
$ws.Cells.Item(73, 9).Value = 4999.524
$ws.Cells.Item(73, 10).Value = 22411.4
$ws.Cells.Item(73, 11).Value = 4999.524
$ws.Cells.Item(73, 12).Value = 22411.4
$ws.Cells.Item(73, 13).Value = -4063.524
$ws.Cells.Item(73, 14).Value = -24283.4
$ws.Cells.Item(102, 8).Value = 50002004
$ws.Cells.Item(102, 9).Value = 1598.3
$ws.Cells.Item(102, 11).Value = 1598.3
$ws.Cells.Item(102, 13).Value = 23.70000000000005
$ws.Cells.Item(122, 8).Value = 3199.875
$ws.Cells.Item(122, 9).Value = 1950
$ws.Cells.Item(122, 11).Value = 5850
$ws.Cells.Item(122, 13).Value = -3400
$ws.Cells.Item(126, 8).Value = 9452.388999999999
$ws.Cells.Item(126, 9).Value = 13278
$ws.Cells.Item(126, 11).Value = 39834
$ws.Cells.Item(126, 13).Value = -37364
$ws.Cells.Item(132, 8).Value = 5043
$ws.Cells.Item(132, 9).Value = 4360.2
$ws.Cells.Item(132, 11).Value = 13080.6
$ws.Cells.Item(132, 13).Value = -10550.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 6443.273
$ws.Cells.Item(40, 9).Value = 4312.1665
$ws.Cells.Item(40, 11).Value = 4312.1665
$ws.Cells.Item(40, 13).Value = -4176.1665
$ws.Cells.Item(46, 8).Value = 1443.2174
$ws.Cells.Item(46, 9).Value = 1260.8
$ws.Cells.Item(46, 10).Value = 1785.25
$ws.Cells.Item(46, 11).Value = 1260.8
$ws.Cells.Item(46, 12).Value = 1785.25
$ws.Cells.Item(46, 13).Value = -1072.8
$ws.Cells.Item(46, 14).Value = -2161.25
$ws.Cells.Item(68, 8).Value = 2986.25
$ws.Cells.Item(68, 9).Value = 2986.25
$ws.Cells.Item(68, 11).Value = 2986.25
$ws.Cells.Item(68, 13).Value = -2237.25
$ws.Cells.Item(71, 8).Value = 2986.25
$ws.Cells.Item(71, 9).Value = 2986.25
$ws.Cells.Item(71, 11).Value = 14931.25
$ws.Cells.Item(71, 13).Value = -11187.25
$ws.Cells.Item(82, 8).Value = 2114.8667
$ws.Cells.Item(82, 9).Value = 2087.4285
$ws.Cells.Item(82, 10).Value = 2499
$ws.Cells.Item(82, 11).Value = 2087.4285
$ws.Cells.Item(82, 12).Value = 2499
$ws.Cells.Item(82, 13).Value = -1726.4285
$ws.Cells.Item(82, 14).Value = -3221
$ws.Cells.Item(85, 8).Value = 2114.8667
$ws.Cells.Item(85, 9).Value = 2087.4285
$ws.Cells.Item(85, 10).Value = 2499
$ws.Cells.Item(85, 11).Value = 2087.4285
$ws.Cells.Item(85, 12).Value = 2499
$ws.Cells.Item(85, 13).Value = -839.4285
$ws.Cells.Item(85, 14).Value = -4995
$ws.Cells.Item(122, 8).Value = 4552.0938
$ws.Cells.Item(122, 9).Value = 3894.25
$ws.Cells.Item(122, 10).Value = 5648.5
$ws.Cells.Item(122, 11).Value = 11682.75
$ws.Cells.Item(122, 12).Value = 16945.5
$ws.Cells.Item(122, 13).Value = -9232.75
$ws.Cells.Item(122, 14).Value = -21845.5
$ws.Cells.Item(131, 8).Value = 32333.334
$ws.Cells.Item(131, 9).Value = 44000
$ws.Cells.Item(131, 11).Value = 44000
$ws.Cells.Item(131, 13).Value = -38960

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 2904.9546
$ws.Cells.Item(136, 10).Value = 2623.0833
$ws.Cells.Item(136, 12).Value = 7869.249899999999
$ws.Cells.Item(136, 14).Value = -12969.2499
